$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.338.06'
$ws.Cells.Item(2, 5).Value = '  -5.17%  '

$ws.Cells.Item(3, 4).Value = '3.350.60'
$ws.Cells.Item(3, 5).Value = '  -4.94%  '

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).Value = '554.74'
$ws.Cells.Item(5, 5).Value = '  -1.09%  '

$ws.Cells.Item(6, 4).Value = '169.64'
$ws.Cells.Item(6, 5).Value = '  -9.67%  '

$ws.Cells.Item(7, 4).Value = '0.605'
$ws.Cells.Item(7, 5).Value = '  -2.51%  '

$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.09%  '

$ws.Cells.Item(9, 4).Value = '0.606'
$ws.Cells.Item(9, 5).Value = '  -4.05%  '

$ws.Cells.Item(10, 5).Value = '  -2.26%  '

$ws.Cells.Item(11, 4).Value = '54.80'
$ws.Cells.Item(11, 5).Value = '  +1.17%  '

$ws.Cells.Item(12, 4).Value = '0.0000262'
$ws.Cells.Item(12, 5).Value = '  -2.84%  '

$ws.Cells.Item(13, 4).Value = '8.83'
$ws.Cells.Item(13, 5).Value = '  -5.13%  '

$ws.Cells.Item(14, 4).Value = '3.880.21'
$ws.Cells.Item(14, 5).Value = '  -4.98%  '

$ws.Cells.Item(15, 4).Value = '0.117'
$ws.Cells.Item(15, 5).Value = '  -3.47%  '

$ws.Cells.Item(16, 4).Value = '3.316.46'
$ws.Cells.Item(16, 5).Value = '  -5.80%  '

$ws.Cells.Item(17, 4).Value = '17.60'
$ws.Cells.Item(17, 5).Value = '  -4.65%  '

$ws.Cells.Item(18, 4).Value = '63.150.72'
$ws.Cells.Item(18, 5).Value = '  -5.39%  '

$ws.Cells.Item(19, 4).Value = '11.50'
$ws.Cells.Item(19, 5).Value = '  -4.14%  '

$ws.Cells.Item(20, 4).Value = '0.966'
$ws.Cells.Item(20, 5).Value = '  -3.03%  '

$ws.Cells.Item(21, 4).Value = '400.70'
$ws.Cells.Item(21, 5).Value = '  -5.61%  '

$ws.Cells.Item(22, 5).Value = '  -1.72%  '

$ws.Cells.Item(23, 5).Value = '  +4.43%  '

$ws.Cells.Item(24, 4).Value = '81.58'
$ws.Cells.Item(24, 5).Value = '  -4.57%  '

$ws.Cells.Item(25, 4).Value = '12.96'
$ws.Cells.Item(25, 5).Value = '  +5.56%  '

$ws.Cells.Item(26, 4).Value = '10.67'
$ws.Cells.Item(26, 5).Value = '  -3.27%  '

$ws.Cells.Item(27, 4).Value = '2.71'
$ws.Cells.Item(27, 5).Value = '  -7.23%  '

$ws.Cells.Item(28, 4).Value = '8.55'
$ws.Cells.Item(28, 5).Value = '  -5.00%  '

$ws.Cells.Item(29, 4).Value = '29.01'
$ws.Cells.Item(29, 5).Value = '  -4.43%  '

$ws.Cells.Item(30, 4).Value = '6.49'
$ws.Cells.Item(30, 5).Value = '  -2.15%  '

$ws.Cells.Item(31, 4).Value = '584.65'
$ws.Cells.Item(31, 5).Value = '  -6.84%  '

$ws.Cells.Item(32, 4).Value = '11.24'
$ws.Cells.Item(32, 5).Value = '  -4.08%  '

$ws.Cells.Item(33, 5).Value = '  -5.63%  '

$ws.Cells.Item(34, 4).Value = '57.79'
$ws.Cells.Item(34, 5).Value = '  -3.74%  '

$ws.Cells.Item(35, 4).Value = '0.149'
$ws.Cells.Item(35, 5).Value = '  +0.44%  '

$ws.Cells.Item(36, 4).Value = '1.00'
$ws.Cells.Item(36, 5).Value = '  +0.13%  '

$ws.Cells.Item(37, 4).Value = '35.60'
$ws.Cells.Item(37, 5).Value = '  -6.78%  '

$ws.Cells.Item(38, 4).Value = '3.42'
$ws.Cells.Item(38, 5).Value = '  +2.29%  '

$ws.Cells.Item(39, 4).Value = '3.145.99'
$ws.Cells.Item(39, 5).Value = '  +0.74%  '

$ws.Cells.Item(40, 4).Value = '0.364'
$ws.Cells.Item(40, 5).Value = '  -5.63%  '

$ws.Cells.Item(41, 5).Value = '  -11.66%  '

$ws.Cells.Item(42, 4).Value = '0.996'
$ws.Cells.Item(42, 5).Value = '  -0.30%  '

$ws.Cells.Item(43, 4).Value = '2.79'
$ws.Cells.Item(43, 5).Value = '  -2.53%  '

$ws.Cells.Item(44, 5).Value = '  -5.97%  '

$ws.Cells.Item(45, 4).Value = '3.19'
$ws.Cells.Item(45, 5).Value = '  -5.20%  '

$ws.Cells.Item(46, 4).Value = '0.0398'
$ws.Cells.Item(46, 5).Value = '  -4.30%  '

$ws.Cells.Item(47, 4).Value = '2.58'
$ws.Cells.Item(47, 5).Value = '  -6.07%  '

$ws.Cells.Item(48, 5).Value = '  -3.70%  '

$ws.Cells.Item(49, 4).Value = '133.37'
$ws.Cells.Item(49, 5).Value = '  -5.34%  '

$ws.Cells.Item(50, 4).Value = '8.05'

$ws.Cells.Item(51, 4).Value = '2.81'
$ws.Cells.Item(51, 5).Value = '  +0.15%  '
